$d = $word.ActiveDocument
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Paragraph 1: "Nguyen huu tinh" -> split into runs with spell-check proofErr marks ---
$p1 = $d.Paragraphs(1).Range
$frag1 = '<w:p ' + $w + '>' +
    '<w:r><w:t xml:space="preserve">Nguyen </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>huu</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>tinh</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '</w:p>'
$p1.InsertXML($frag1)

# --- Paragraph 2: "Nguyen Quang Vinh" -> split with proofErr marks, plus a new line
#     "Pham Minh Taii " added via a line break, keeping the trailing bookmark ---
$p2 = $d.Paragraphs(2).Range
$frag2 = '<w:p ' + $w + '>' +
    '<w:r><w:t xml:space="preserve">Nguyen </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Quang</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Vinh</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:br/><w:t xml:space="preserve">Pham Minh </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Taii</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>'
$p2.InsertXML($frag2)
